$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header text updates (new report week / volume number)
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 30   Number  44"
$ws.Range("C9").Value = "Report Covering the Week  10/30/2023  Through  11/5/2023"

# ---------------------------------------------------------------------------
# Helper: cells that flip from a NUMBER to the literal text "0" (style 14,
# General format, shared string "0") need the quote-prefix trick so the
# numeric-looking text isn't silently re-coerced back into a number, then a
# PasteSpecial(xlPasteFormats) from an existing style-14 cell to pick up the
# exact right-aligned / General-format style without minting a new style.
# Cells that flip to "***.*" don't need the quote prefix since that text can
# never be parsed as a number. Cells that flip from text back to a plain
# number just need PasteSpecial(xlPasteFormats) from an existing numeric
# (style 15) cell followed by the numeric value.
# ---------------------------------------------------------------------------

# Row 15 (Rape): G15 1 -> "0" ; H15 -100 -> "***.*"
$ws.Range("G15").Value = "'0"
$ws.Range("F15").Copy()
$ws.Range("G15").PasteSpecial(-4122)

$ws.Range("H15").Value = "***.*"
$ws.Range("F15").Copy()
$ws.Range("H15").PasteSpecial(-4122)

# Row 16 (Robbery)
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 200
$ws.Range("F16").Value = 22
$ws.Range("G16").Value = 12
$ws.Range("H16").Value = 83.333333333333
$ws.Range("I16").Value = 199
$ws.Range("J16").Value = 170
$ws.Range("K16").Value = 17.058823529411
$ws.Range("L16").Value = 74.561403508771
$ws.Range("M16").Value = -11.160714285714
$ws.Range("N16").Value = -76.393831553973

# Row 17 (Fel. Assault)
$ws.Range("C17").Value = 9
$ws.Range("D17").Value = 7
$ws.Range("E17").Value = 28.571428571428
$ws.Range("F17").Value = 25
$ws.Range("G17").Value = 28
$ws.Range("H17").Value = -10.714285714285
$ws.Range("I17").Value = 242
$ws.Range("J17").Value = 259
$ws.Range("K17").Value = -6.563706563706
$ws.Range("L17").Value = 10.502283105022
$ws.Range("M17").Value = 29.411764705882
$ws.Range("N17").Value = -7.279693486590

# Row 18 (Burglary)
$ws.Range("C18").Value = 5
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 150
$ws.Range("F18").Value = 18
$ws.Range("G18").Value = 18
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 173
$ws.Range("J18").Value = 233
$ws.Range("K18").Value = -25.751072961373
$ws.Range("L18").Value = -7.978723404255
$ws.Range("M18").Value = -55.064935064935
$ws.Range("N18").Value = -89.859320046893

# Row 19 (Gr. Larceny)
$ws.Range("C19").Value = 17
$ws.Range("D19").Value = 16
$ws.Range("E19").Value = 6.25
$ws.Range("F19").Value = 65
$ws.Range("G19").Value = 60
$ws.Range("H19").Value = 8.333333333333
$ws.Range("I19").Value = 594
$ws.Range("J19").Value = 559
$ws.Range("K19").Value = 6.261180679785
$ws.Range("L19").Value = 30.837004405286
$ws.Range("M19").Value = 61.413043478260
$ws.Range("N19").Value = 9.191176470588

# Row 20 (G.L.A.)
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = 100
$ws.Range("F20").Value = 22
$ws.Range("G20").Value = 28
$ws.Range("H20").Value = -21.428571428571
$ws.Range("I20").Value = 322
$ws.Range("J20").Value = 258
$ws.Range("K20").Value = 24.806201550387
$ws.Range("L20").Value = 81.920903954802
$ws.Range("M20").Value = 2.222222222222
$ws.Range("N20").Value = -89.27381745503

# Row 21 (TOTAL)
$ws.Range("C21").Value = 40
$ws.Range("D21").Value = 29
$ws.Range("E21").Value = 37.931034482758
$ws.Range("F21").Value = 153
$ws.Range("G21").Value = 146
$ws.Range("H21").Value = 4.794520547945
$ws.Range("I21").Value = 1553
$ws.Range("J21").Value = 1500
$ws.Range("K21").Value = 3.533333333333
$ws.Range("L21").Value = 33.190394511149
$ws.Range("M21").Value = 3.533333333333
$ws.Range("N21").Value = -75.707805412169

# Row 22 (Transit): C22 "0" -> 1 (text to number)
$ws.Range("I16").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("C22").Value = 1
$ws.Range("F22").Value = 2
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = 100
$ws.Range("I22").Value = 20
$ws.Range("K22").Value = 53.846153846153
$ws.Range("L22").Value = 233.333333333333
$ws.Range("M22").Value = 11.111111111111

# Row 23 (Housing) - unchanged

# Row 24 (Petit Larceny)
$ws.Range("C24").Value = 16
$ws.Range("D24").Value = 23
$ws.Range("E24").Value = -30.434782608695
$ws.Range("F24").Value = 98
$ws.Range("G24").Value = 118
$ws.Range("H24").Value = -16.949152542372
$ws.Range("I24").Value = 1121
$ws.Range("J24").Value = 1245
$ws.Range("K24").Value = -9.959839357429
$ws.Range("L24").Value = -0.971731448763
$ws.Range("M24").Value = 24.972129319955

# Row 25 (Misd. Assault)
$ws.Range("C25").Value = 17
$ws.Range("D25").Value = 18
$ws.Range("E25").Value = -5.555555555555
$ws.Range("F25").Value = 45
$ws.Range("G25").Value = 46
$ws.Range("H25").Value = -2.173913043478
$ws.Range("I25").Value = 409
$ws.Range("J25").Value = 463
$ws.Range("K25").Value = -11.663066954643
$ws.Range("L25").Value = -4.215456674473
$ws.Range("M25").Value = -33.387622149837

# Row 26 (UCR Rape*): F26 1 -> "0" (number to text)
$ws.Range("F26").Value = "'0"
$ws.Range("F15").Copy()
$ws.Range("F26").PasteSpecial(-4122)

$ws.Range("H26").Value = -100
$ws.Range("J26").Value = 29
$ws.Range("K26").Value = 3.448275862068

# Row 27 (Other Sex Crimes): C27 "0" -> 1 ; D27 3 -> "0" ; E27 -100 -> "***.*"
$ws.Range("I16").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("C27").Value = 1

$ws.Range("D27").Value = "'0"
$ws.Range("F15").Copy()
$ws.Range("D27").PasteSpecial(-4122)

$ws.Range("E27").Value = "***.*"
$ws.Range("F15").Copy()
$ws.Range("E27").PasteSpecial(-4122)

$ws.Range("F27").Value = 3
$ws.Range("G27").Value = 6
$ws.Range("I27").Value = 54
$ws.Range("K27").Value = -14.285714285714
$ws.Range("L27").Value = 20
